# "add glob finder and baseCommand check"
#
# Adds a second data row to the "Tabelle1" annotation table, mirroring the
# existing row but pointing at a second glob/file pair (hello2.tar /
# hello2.txt), grows the Excel table to cover the new row, renames the
# built-in cell style from its localized "Standard" name to "Normal", and
# leaves the "Tabelle1" worksheet as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Add new row 3 -------------------------------------------------------
# Column order matters for shared-string allocation: set K3 ("Derived Data
# File" -> hello2.txt) before A3 ("Source Name" -> hello2.tar) so new shared
# strings are appended in the same order Excel produced them.
$ws.Range("K3").Value = "hello2.txt"
$ws.Range("A3").Value = "hello2.tar"
$ws.Range("B3").Value = 0

# E3 ("Characteristic [Prefix]") reuses the existing "--file" shared string
# and the text-format style applied to E2; copy E2 wholesale so the value and
# style land together.
$ws.Range("E2").Copy($ws.Range("E3"))

$ws.Range("H3").Value = "tar --extract"

# --- Grow the Excel table to include the new row -------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K3"))

# --- Rename the built-in cell style "Standard" -> "Normal" ---------------
# The workbook's only cell style is the built-in Normal style, saved under
# its localized name "Standard". Excel never really deletes the built-in
# Normal style; requesting its deletion instead resets it back to its
# canonical "Normal" name.
$wb.Styles.Item(1).Delete()

# --- Make "Tabelle1" the active sheet / selection -------------------------
$ws.Activate()
$ws.Range("B8").Select()
